$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '56.492.41'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.10%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.309.34'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.37%  '

$ws.Range('E4').Value = '  -0.06%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '517.94'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.51%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '131.49'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.50%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.996'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.20%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.532'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.72%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.331.66'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.11%  '

$ws.Range('E10').Value = '  -1.73%  '

$ws.Range('E11').Value = '  +0.06%  '

$ws.Range('E12').Value = '  -1.43%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.338'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.35%  '

$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.728.08'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.14%  '

$ws.Range('B15').Value = 'Avalanche'
$ws.Range('C15').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '23.33'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.06%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '56.488.64'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.13%  '

$ws.Range('E17').Value = '  -1.60%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.322.53'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.64%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '333.25'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.37%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.35'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.31%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.13'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.90%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.72'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.26%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.997'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.17%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '61.15'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.90%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '8.64'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +8.76%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.164'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.21%  '

$ws.Range('E27').Value = '  -0.37%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.32'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.80%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '169.62'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.24%  '

$ws.Range('E30').Value = '  -0.15%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0₃0714'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.22%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.09'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.33%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '18.29'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.95%  '

$ws.Range('E34').Value = '  -0.04%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.996'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.19%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.24'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.77%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.91'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.89%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.881'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.59%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.57'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.22%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '38.77'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.36%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '148.35'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +5.07%  '

$ws.Range('E42').Value = '  -1.52%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '286.04'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.49%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.57'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.75%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '5.05'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.07%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0925'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.02%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0498'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.31%  '

$ws.Range('E48').Value = '  -0.21%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '18.35'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +3.13%  '

$ws.Range('B50').Value = 'Polygon'
$ws.Range('C50').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.377'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.64%  '

$ws.Range('B51').Value = 'VeChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0213'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.02%  '
